$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("subject")

# Row 2 updates: subject_label becomes a plain number, source_study
# changes from "ARC" to "IAVI", and the rest of the row shifts as the
# shared-string table is rewritten (handled implicitly by setting values).
$ws.Range("A2").Value = 175055
$ws.Range("B2").Value = "IAVI"

# Update the visible selection to B3 (matches the new <selection> element).
$ws.Range("B3").Select()
